$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 7623.3335
$ws.Range("I19").Value = 25562.75
$ws.Range("J19").Value = 1099.909
$ws.Range("K19").Value = 25562.75
$ws.Range("L19").Value = 1099.909
$ws.Range("M19").Value = -25387.75
$ws.Range("N19").Value = -1449.909
$ws.Range("H111").Value = 3566.8572
$ws.Range("I111").Value = 2019
$ws.Range("J111").Value = 3824.8333
$ws.Range("K111").Value = 6057
$ws.Range("L111").Value = 11474.4999
$ws.Range("M111").Value = -2990
$ws.Range("N111").Value = -17608.4999
$ws.Range("H137").Value = 1583.871
$ws.Range("I137").Value = 1272.9615
$ws.Range("J137").Value = 3200.6
$ws.Range("K137").Value = 3818.8845
$ws.Range("L137").Value = 9601.799999999999
$ws.Range("M137").Value = -1268.8845
$ws.Range("N137").Value = -14701.8

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4989.3335
$ws.Range("I32").Value = 2667.7737
$ws.Range("K32").Value = 2667.7737
$ws.Range("M32").Value = -2380.7737

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 619.6667
$ws.Range("J64").Value = 662.8333
$ws.Range("L64").Value = 662.8333
$ws.Range("N64").Value = -1112.8333
$ws.Range("H67").Value = 619.6667
$ws.Range("J67").Value = 662.8333
$ws.Range("L67").Value = 662.8333
$ws.Range("N67").Value = -2222.8333

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1762.7142
$ws.Range("I16").Value = 1882
$ws.Range("J16").Value = 1325.3334
$ws.Range("K16").Value = 1882
$ws.Range("L16").Value = 1325.3334
$ws.Range("M16").Value = -1595
$ws.Range("N16").Value = -1899.3334
$ws.Range("H31").Value = 2026.9535
$ws.Range("I31").Value = 1671.9667
$ws.Range("K31").Value = 1671.9667
$ws.Range("M31").Value = -1376.9667
$ws.Range("H34").Value = 2026.9535
$ws.Range("I34").Value = 1671.9667
$ws.Range("K34").Value = 1671.9667
$ws.Range("M34").Value = -1469.9667
$ws.Range("H58").Value = 1620.878
$ws.Range("I58").Value = 1112.3182
$ws.Range("J58").Value = 2209.7368
$ws.Range("K58").Value = 1112.3182
$ws.Range("L58").Value = 2209.7368
$ws.Range("M58").Value = -909.3181999999999
$ws.Range("N58").Value = -2615.7368
$ws.Range("H99").Value = 100000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 100000
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 100000
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -102996
$ws.Range("H113").Value = 1762.7142
$ws.Range("I113").Value = 1882
$ws.Range("J113").Value = 1325.3334
$ws.Range("K113").Value = 1882
$ws.Range("L113").Value = 1325.3334
$ws.Range("M113").Value = 288
$ws.Range("N113").Value = -5665.3334
$ws.Range("H126").Value = 100000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 100000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 300000
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -304940
$ws.Range("H132").Value = 2806.963
$ws.Range("I132").Value = 1920.2222
$ws.Range("J132").Value = 4580.4443
$ws.Range("K132").Value = 5760.6666
$ws.Range("L132").Value = 13741.3329
$ws.Range("M132").Value = -3230.6666
$ws.Range("N132").Value = -18801.3329
$ws.Range("H134").Value = 2104.8386
$ws.Range("I134").Value = 1981.6154
$ws.Range("J134").Value = 2745.6
$ws.Range("K134").Value = 5944.8462
$ws.Range("L134").Value = 8236.799999999999
$ws.Range("M134").Value = -3409.8462
$ws.Range("N134").Value = -13306.8
$ws.Range("H136").Value = 1620.878
$ws.Range("I136").Value = 1112.3182
$ws.Range("J136").Value = 2209.7368
$ws.Range("K136").Value = 3336.9546
$ws.Range("L136").Value = 6629.2104
$ws.Range("M136").Value = -786.9546
$ws.Range("N136").Value = -11729.2104

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 196.42105
$ws.Range("I14").Value = 196.42105
$ws.Range("K14").Value = 589.26315
$ws.Range("M14").Value = -416.26315
$ws.Range("H33").Value = 158.65218
$ws.Range("I33").Value = 116.9375
$ws.Range("J33").Value = 254
$ws.Range("K33").Value = 701.625
$ws.Range("L33").Value = 1524
$ws.Range("M33").Value = -418.625
$ws.Range("N33").Value = -2090
$ws.Range("H44").Value = 1137.875
$ws.Range("I44").Value = 434.33334
$ws.Range("J44").Value = 1560
$ws.Range("K44").Value = 1303.00002
$ws.Range("L44").Value = 4680
$ws.Range("M44").Value = -905.0000199999999
$ws.Range("N44").Value = -5476
$ws.Range("H113").Value = 2086.7778
$ws.Range("I113").Value = 2600
$ws.Range("J113").Value = 1573.5555
$ws.Range("K113").Value = 7800
$ws.Range("L113").Value = 4720.666499999999
$ws.Range("M113").Value = -5630
$ws.Range("N113").Value = -9060.666499999999
$ws.Range("H122").Value = 665
$ws.Range("I122").Value = 368
$ws.Range("J122").Value = 720.6875
$ws.Range("K122").Value = 3312
$ws.Range("L122").Value = 6486.1875
$ws.Range("M122").Value = -862
$ws.Range("N122").Value = -11386.1875
$ws.Range("H131").Value = 2892
$ws.Range("I131").Value = 422.72726
$ws.Range("J131").Value = 3509.318
$ws.Range("K131").Value = 1268.18178
$ws.Range("L131").Value = 10527.954
$ws.Range("M131").Value = 3771.81822
$ws.Range("N131").Value = -20607.954
$ws.Range("H132").Value = 817.0833
$ws.Range("I132").Value = 800
$ws.Range("J132").Value = 868.3333
$ws.Range("K132").Value = 7200
$ws.Range("L132").Value = 7814.9997
$ws.Range("M132").Value = -4670
$ws.Range("N132").Value = -12874.9997

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1245.1875
$ws.Range("I113").Value = 1232.3077
$ws.Range("J113").Value = 1301
$ws.Range("K113").Value = 1232.3077
$ws.Range("L113").Value = 1301
$ws.Range("M113").Value = 937.6922999999999
$ws.Range("N113").Value = -5641
$ws.Range("H132").Value = 3060.5293
$ws.Range("I132").Value = 2328.25
$ws.Range("J132").Value = 3711.4443
$ws.Range("K132").Value = 6984.75
$ws.Range("L132").Value = 11134.3329
$ws.Range("M132").Value = -4454.75
$ws.Range("N132").Value = -16194.3329

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 47624548
$ws.Range("I136").Value = 6500
$ws.Range("J136").Value = 66671770
$ws.Range("K136").Value = 19500
$ws.Range("L136").Value = 200015310
$ws.Range("M136").Value = -16950
$ws.Range("N136").Value = -200020410

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2663.25
$ws.Range("I132").Value = 2357.8572
$ws.Range("J132").Value = 3090.8
$ws.Range("K132").Value = 7073.571599999999
$ws.Range("L132").Value = 9272.400000000001
$ws.Range("M132").Value = -4543.571599999999
$ws.Range("N132").Value = -14332.4
$ws.Range("H136").Value = 31241.184
$ws.Range("I136").Value = 46138.184
$ws.Range("J136").Value = 10757.8125
$ws.Range("K136").Value = 138414.552
$ws.Range("L136").Value = 32273.4375
$ws.Range("M136").Value = -135864.552
$ws.Range("N136").Value = -37373.4375

$wb.Save()
Write-Host "Applied all Carbuncle_Profits cell updates."